$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Nike Dunk Low ''Panda Black White''', 135),
    @('Nike Air Force 1 Low ''07 White''', 150),
    @('Maison MIHARA YASUHIRO PETERSON OG Sole Canvas Low ''Black''', 125),
    @('Nike Air Max 95 SP ''Corteiz Pink Beam''', 130),
    @('Samba OG ''Cloud White Core Black''', 130),
    @('Nike Air Huarache ultra ''Triple Black''', 125),
    @('Nike Air Huarache ultra ''triple white''', 125),
    @('Maison Mihara Yasuhiro Peterson  Low ''Black Overhanging'' ', 125),
    @('Maison MIHARA YASUHIRO Peterson x Amiri ''White''', 125),
    @('Maison MIHARA YASUHIRO Peterson x Amiri ''Black''', 125),
    @('MAISON MIHARA YASUHIRO BLAKEY OG SOLE CANVAS LOW ''White''', 125),
    @('MAISON MIHARA YASUHIRO BLAKEY OG SOLE CANVAS LOW ''Black/White''', 125),
    @('Handball Spezial ''Light Blue''', 130),
    @('Nike SB Dunk Low ''Mummy''', 95),
    @('Birkenstock BLACK', 95),
    @('Maison Mihara Yasuhiro Peterson  Low ''White Overhanging'' ', 125),
    @('Manchester SPZL x C.P. Company ''Light Blue''', 130),
    @('MAISON MIHARA YASUHIRO BLAKEY OG SOLE CANVAS LOW ''Black''', 125),
    @('Saucony Jazz Original Vintage ''Dark Grey''', 150),
    @('DC Shoes Court Graffik ''Black Sabbath''', 95),
    @('Nike SB Dunk Low ''Night of Mischief''', 95),
    @('Yeezy Boost 350 V2 ''Black Non-Reflective''', 120),
    @('Converse X COMME Des GARCONS PLAY Chuck 70 Hi ''Milk''', 65),
    @('Globe Tilt trainers ''Strom Grey/Black''', 150),
    @('Nike SB Dunk Low ''Club 58 Gulf''', 95),
    @('Handball Spezial ''Black Gum''', 130),
    @('Nike Air Zoom Pegasus 39 ''Black Anthracite''', 120),
    @('Nike Air Jordan 1 Low ''Court Purple''', 95),
    @('Nike Air Zoom Vomero 5 x off-white ''Supersonic''', 120),
    @('Saucony ProGrid Omni 9 ''White blue''', 95),
    @('ASICS Gel-Kayano 14 ''Illuminate Yellow''', 140),
    @('ASICS Gel-Kayano 14 ''Clay Grey Pure Silver''', 95),
    @('Maison MIHARA YASUHIRO PETERSON NEW Black RAF', 125),
    @('Maison MIHARA YASUHIRO PETERSON NEW White/Black', 125),
    @('Climacool 1 ''White/Black/Red''', 125),
    @('Maison MIHARA YASUHIRO PETERSON OG Sole Canvas Low ''All Black''', 125),
    @('Song for the Mute Adistar Cushion ''Taupe Metallic Granite''', 150),
    @('Nike Air Portal x SUPREME ''Black/Brown''', 120),
    @('Converse Canvas High ''Tawny Owl''', 65),
    @('Nike SB Dunk Low Tiffany & Co', 95),
    @('Nike ACG Mountain Fly GTX ''Clay Green''', 85),
    @(' ASICS Gel-Kahana TR V4 ''Carbon Black''', 140),
    @('Handball Spezial ''Grey White/ Gum''', 130),
    @('Nike v2k run ''Summit White/Metallic Silver''', 140),
    @('Handball SPEZIAL ''ACTIVE GREEN GUM''', 130),
    @('DC Stag White/grey/blue', 150),
    @('DC Stag ''Black/Brown''', 150),
    @('Saucony Jazz Original Vintage ''Black/Gym''', 150),
    @('DC Shoes Court Graffik ''Marvel''', 95),
    @('DC Shoes Court Graffik ''Black/White''', 95),
    @('Asics GEL Kahana 8 Silver/Ruby', 140),
    @('Samba Team ''Argentina''', 130),
    @('Samba OG ''Black White Gum''', 130),
    @('Nike ACG Mountain Fly Gore-Tex ''Khaki''', 85),
    @('Forum Buckle Low x Bad Bunny ''Last Forum''', 95),
    @('Yeezy 350 v2 ''White''', 120),
    @('Nike Air Zoom Vomero 5 x off-white ''Superkevin''', 120),
    @('XLG Metawave ', 140),
    @('Nike SB Dunk x Chrome Hearts Brown', 95),
    @('ClimaCool 1 ''Noir''', 125),
    @('Nike Air Max Pulse Roam ''Cool Grey''', 130),
    @('Nike React Vision ''Double Swoosh''', 135),
    @('Nike Air Max TL 2.5 ''Grey''', 150),
    @('Nike X Supreme Shox Ride 2 SP ''White/Navy/Red''', 130),
    @('Nike Free Run 2 ''Black/White''', 125),
    @('Nike Free Run 2 ''Triple Black''', 125),
    @('Nike React Element 87 ''Anthracite''', 135),
    @('Nike React Element 55 ''Triple Black''', 135),
    @('Nike React Element 55 ''Reverse Schematic - Black''', 135),
    @('Nike React Element 87 ''Solar Red''', 135),
    @('Nike SB Dunk Low Racing Buggy', 95),
    @('Nike SB Dunk Low x NBA Kobe Bryant', 95),
    @('Nike SB Dunk Low Travis Scott x Play Station 5 ''Khakki''', 95),
    @('DC Shoes Court Graffik ''White & Black''', 95),
    @('Nike Air Max 95 Sneaker Boot ''Triple Black''', 85),
    @('Nike ACG Mountain Fly GTX ''Black''', 100),
    @('DVS Comanche Grey/White/Yellow/Black', 150),
    @('Nike Air Jordan 1 mid x Kaws ''Black''', 95),
    @('Nike Dunk Low Pro SB ''Purple Pigeon''', 95),
    @('Nike Air Max 95 x CORTEIZ ''Gutta Green'' ', 130),
    @('Nike React Vision ''Smoke Grey Orange''', 135),
    @('Nike React Vision ''Triple Black''', 135),
    @('DC Shoes Manteca 4 x ''Slayer''', 95),
    @('Nike P-6000 ''Athletic Department''', 95),
    @('Superstar ''Core Black/White''', 135),
    @('DC Stag ''black/orange''', 150),
    @('DC Stag ''Black/White''', 150),
    @('DC Stag ''Black/Red''', 150),
    @('NIKE Air Force 1 Wild ''black''', 75),
    @('ASICS Gel-Kahana 8 ''Black''', 150),
    @('Converse Run Star Motion ''Black Gum Honey''', 120),
    @('Forum 84 Low Shoes ''Blue Denim Gum''', 95),
    @('Converse X COMME Des GARCONS PLAY Chuck 70 Hi ''Black Red''', 65),
    @('Nike Air Max 90 Essential ''Hyper Jade''', 150),
    @('Yeezy Boost 350 V2 ''Desert Sage''', 120),
    @('yeezy 350 v2 grey orange', 120),
    @('Nike Air Zoom Pegasus 39 ''Black White''', 120),
    @('Nike Air Zoom Pegasus 39 TB ''White Black''', 120),
    @('Birkenstock Chocolate', 95),
    @('Nike React Vision ''Wolf Grey''', 135),
    @('Nike React Vision ''White''', 135),
    @('ASICS Gel-Kayano 14 ''Ironclad Cloud Grey''', 150),
    @('ASICS Gel Quantum Kinetic ''Dark Sepia Shamrock''', 150),
    @('Nike Air Max 95 ''Neutral Olive''', 130),
    @('Nike Air Max 95 Essential ''Halloween''', 130),
    @('Nike Air Max Muse ''Triple Black''', 120),
    @('Nike X Concepts SB Dunk Low ''Green Lobster''', 95),
    @('Nike React Element 55 ''Black White''', 135),
    @('Nike SB Dunk Low Travis Scott x Play Station 5 ''grey''', 95),
    @('Nike Dunk Low Batman & Joker', 95),
    @('DC Shoes Court Graffik ''Black Gum''', 95),
    @('Asics Gel-Kahana 8 ''Anthracite'' Gore-Tex', 150),
    @('ASICS Gel-Kahana 8 ''Retro Grey'' Gore-tex ', 150),
    @('Salomon XT-6 GORE-TEX  Series ''Black/Lunar Rock''  ', 95),
    @('DC SHOES SLAYER PURE', 95),
    @('Salomon XT-6 Series ''Night sky/ Blue lime''', 95),
    @('Oztral ''Black/Purple''', 150),
    @('Nike Air Max 96 ''Mystic Teal''', 140),
    @('Nike Air Max 95 X Corteiz ''Aegean Storm''', 130),
    @('Yeezy Boost 350 V2 ''Sand Taupe''', 120),
    @('Ozmillen ''Grey Shadow Violet''', 120),
    @('DC SHOES COURT GRAFFIK X ''SLAYER''', 95),
    @('Nike ACG Air Terra Antarktik Gore-Tex ''Anthracite''', 85),
    @('Nike ACG Air Terra Antarktik Gore-Tex ''Juniper Fog''', 85),
    @('Nike ACG Terra Antarktik GTX ''Black''', 85),
    @('Ocai Retro x Ann Andelman ''Jeans''', 110),
    @('Campus 00S ''Energy Ink''', 75),
    @('Nike Zoom Vomero 5 ''Black Sesame''', 140),
    @('Ocai 001 x Ann Andelman ''Denim Black''', 110),
    @('ASICS Gel-Kahana 8 ''Green/Black''', 150),
    @('Nike SB Dunk Low Pro ''Neckface''', 140),
    @('Nike X Grateful Dead SB Dunk Low ''Yellow Bear''', 95),
    @('Forum low x bad bunny ''back to school''', 95),
    @('Nike Blazer Mid ''77 ''Black Sail''', 95),
    @('Nike Travis Scott X Dunk Low Premium QS SB ''Cactus Jack''', 95),
    @('Nike Dunk Low SE ''Sail Neptune Green''', 95),
    @('Nike Air Max 95 OG ''Neon''', 130),
    @('Yeezy Boost 350 V2 ''Lundmark Non-Reflective''', 120),
    @('Nike Air Max 95/97 x Ducks of a Feather ''Lighting''', 130),
    @('Nike Air Max 90 ''University Gold''', 150),
    @('Nike Air Max 95 Essential ''Armory Navy''', 130),
    @('Nike x Stash Air Max 95 ''Obsidian Black''', 130),
    @('Campus X Bad Bunny ''Grey''', 95),
    @('Nike Shox TL ''Sunrise Gradient''', 130),
    @('ASICS Gel-Nimbus 10.1 ''Black Grand Shark''', 150),
    @('SL 72 OG ''Preloved Brown Almost Yellow''', 75),
    @('Nike Air Max 95 ''Black Grape''', 130),
    @('Mizuno Adventure Twe', 140),
    @('Nike Air Max 95 ''Topographic''', 130),
    @('Nike x size? Air Max Dn ''Platinum Violet''', 120),
    @('Nike Air Max Tailwind 4 ''Grey''', 140),
    @('Saucony Jazz Original Vintage ''Black/White''', 150),
    @('Nike React Element 55 ''Tour Yellow''', 135),
    @('Nike Kyrie 7 EP ''Midnight Navy''', 95),
    @('Salomon XT 6 GTX ''Utility Gray''', 95),
    @('Fear of God Athletics Basketball ''Carbon''', 75),
    @('Salomon XT-6 Series GTX ''Blue'' GORE-TEX', 95),
    @('DVS Comanche Black/White', 150),
    @('Air Jordan 1 High "TEATHER Brown&Green& Beige"', 65),
    @('Nike Air Max Plus Drift ''Bright Crimson Field/Purple''', 140),
    @('Nike React Vision ''White Anthracite''', 135),
    @('Nike React Vision ''Black''', 135),
    @('Handball Spezial ''Core Black/Grey Four''', 130),
    @('Nike ZoomX Zegama Trail 2 ''Daybrake Clay''', 120),
    @('Ozmillen ''Grey two/ Crystall white''', 120),
    @('Handball Spezial ''Pantone/Collegiate Navy/Gum''', 130),
    @('Handball Spezial ''Exclusive Mesh Purple''', 130),
    @('Yeezy Boost 350 V2 ''Ash Pearl''', 120),
    @('DC SHOES COURT GRAFFIK ''BASQ CT''', 95),
    @('DC Shoes Court Graffik ''Black/Red''', 95),
    @('Zig Kinetica 2 Edge GTX ''Black Midnight Pine''', 85),
    @('NIKE JORDAN 1 MID x Travis Scott ''BLACK PHANTOM''', 65),
    @('Nike x Off-White Air Rubber Dunk ''Green Strike''', 135),
    @('Samba Team x Real Madrid ''White Gum'' ', 130),
    @('Nike ZoomX Zegama Trail ''Volt''', 120),
    @('Campus 00s ''Black/White''', 75),
    @('Oztral ''Aluminium Grey''', 120),
    @('Handball Spezial ''Navy Gum''', 130),
    @('Asics Gel-Kahana 8 ''Grey/Blue''', 150),
    @('Saucony Layer ''Black''', 65),
    @('ASICS Gel Kayano 14 ''Vanilla Tarmac''', 150),
    @('Nike Air Jordan 4 Retro ''Canyon Purple''', 65),
    @('Nike Air Max 90 Terrascape ''Sail Sea Glass''', 95),
    @('Nike X Grateful Dead SB Dunk Low ''Orange Bear''', 75),
    @('Niteball ''Black Sub Green''', 65),
    @('Nike Air Jordan 4 Retro SE ''Neon 95''', 65),
    @('Superstar ''White/Black''', 135),
    @('Merrell ice cap moc ''black/orange''', 100),
    @('ASICS Gel-NYC ''Mauve Blue''', 140),
    @('Handball Spezial Moscow ', 130),
    @('Nike Air Max 95 ''Koi''', 130),
    @('Nike x Supreme Shox Ride 2 SP White/Grey/Volt', 150),
    @('Nike Air Max 95 Petek ''Dark Blue/Purple''', 130),
    @('Nike Air Max 95 white/black/fiol', 130),
    @('Nike Shox Ride 2 ''Baroque Brown Diffused Blue', 150),
    @('Nike Shox TL ''Persian Violet''', 150),
    @('Nike Nocta Glide Drake ''Khaki''', 140),
    @(' Nike Air Portal x SUPREME ''White/Blue''', 120),
    @('Nike Air Max Portal ''Black/White''', 120),
    @('Nike Air Max 96 ''Triple Black''', 150),
    @('Nike Air Force low ''Venom''', 120),
    @('Nike Dunk Low ''Year of the Dragon''', 95),
    @('Nike Kyrie 7 EP ''Chinese New Year''', 95),
    @('Asics Gel-NYC Gore-Tex ''Black/Bronze''', 95),
    @('Salomon XT-6 GORE-TEX  ''Triple Black''', 95),
    @('ASICS Gel Kayano Trainer 21 ''Black'' ', 140),
    @('SALOMON XT-6 ''Pewter Black Aurora Red'' ', 95),
    @('Nike Dunk Low ''Sanddrift Light Armory Blue''', 95),
    @('Nike Air Jordan 4 retro ''military brown''', 65),
    @('ASICS Gel-Kayano Trainer 21 ''Smoke Blue''', 140),
    @('Nike Air Max Terrascape Plus ''Triple Black''', 95),
    @('Salomon XT-6 Series ''Black/Lunar Rock''', 95),
    @('Globe Tilt trainers ''Dark Shadow/Phantom''', 150),
    @('DC Shoes Spectre x NEEDLES ''Black/Purple''', 95),
    @('Iniki Runner ''Core Black 2.0''', 125),
    @('ASICS Gel 1090 ''Silver''', 140),
    @('Nike Air Zoom Pegasus 38 ''Wolf Grey''', 120),
    @('Nike Zoom Pegasus 39 ''Ashen Slate Football Grey Cobalt Bliss''', 120),
    @('Nike ZoomX Zegama Trail 2 ''Black/White''', 120),
    @('ASICS GEL-1130 X HAL STUDIOS ''Smoke Grey''', 140),
    @('Nike v2k run ''metallic silver/white''', 140),
    @('Asics Gel-1130 ''Slate Grey/Piedmont Grey''', 140),
    @('Asics gel-1130 X HAL STUDIOS  ''Forest''', 140),
    @('Suba Sport black red', 50),
    @('Intimidation Sports Casual Shoes ''Black Purple''', 140),
    @('Premiata Mick ''Black''', 140),
    @('Nike SB Dunk Low Vintage ''Phoenix Suns''', 95),
    @('ASICS GEL SONOMA 15 White/Dark grey', 140),
    @('Saucony Jazz Original Vintage ''Light Grey''', 150),
    @('ASICS Gel-Kahana 8 ''Dark Blue Brown'' demi', 150),
    @('Zig Kinetica 2 Edge GTX ''Black''', 100),
    @('Salomon XT-6 ''Triple Black''', 95),
    @('Nike Air Jordan MVP 678 ''Raptors''', 95),
    @('Nike Air Jordan 6 Retro x Travis Scott ''Olive''', 110),
    @('Nike SB Dunk Low х Stussy ''BMW M3''', 95),
    @('Ocai ''black/smoke grey''', 110),
    @('DC Shoes Court Graffik White', 95),
    @('Campus x bad bunny ''black''', 95),
    @('Nike Air Huarache Craft ''Smokey Grey''', 95),
    @('Nike Air Jordan 1 Low x Travis Scott ''Tiffany''', 75),
    @('Campus 00s ''Ambient Sky''', 75),
    @('Campus 00s ''Grey''', 75),
    @('Nike Zoom Vomero 5 ''Oatmeal''', 140),
    @('Response x Bad Bunny CL ''Cream White''', 140),
    @('A BATHING APE Bape Sta Hi ''black/white''', 65),
    @('ADI2000 ''Lab Green''', 65),
    @('Forum 84 Low ''Off white''', 95),
    @('Nike Dunk Low Pro SB ''Thunderstorm''', 75),
    @('Nike Dunk Low ''Dusty Olive''', 95),
    @('Nike LPL X Dunk Low ''Black White Red''', 75),
    @('Saucony Jazz ''Grey/Black'' New', 95),
    @('Nike Air Max 96 ''Beach''', 140),
    @('Nike SB Dunk Low ''INITIAL D Toyota AE86''', 95),
    @('Nike Lunar Force 1 Duckboot Low ''Ale Brown''', 65),
    @('Nike Air Jordan 4 ''PSG''', 65),
    @('Nike Air Jordan 4 x Travis Scott ''Purple''', 65),
    @('Nike Air Jordan 4 Retro ''Cement''', 65),
    @('Nike Air Jordan 1 Retro High OG ''Crimson Tint''', 65),
    @('Nike Zoom 2K ''White Black''', 95),
    @('Balenciaga Triple S ''White''', 120),
    @('Balenciaga Triple S ''Black''', 120),
    @('Nike Air Monarch IV ''Black''', 95),
    @('Nike Air Force 1 Low ''Triple Black''', 150),
    @('dsquared 2 ''Triple black''', 145),
    @('Nike SF Air Force 1 Mid ''Triple Ivory''', 145),
    @('Nike Air More Uptempo ''Limestone''', 145),
    @('Kappa Fatty ''White/Beige''', 145),
    @('Nike Air Trainer 1 ''Honeydew''', 145),
    @('Nike Air More Uptempo ''Reflective''', 75)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
Write-Host "Updated $($data.Length) rows"
